# Update existing rows 2-4 and add new rows 5-10 for the "ECs" cluster,
# per Dr Hou's advice (adds ECs as sending/target cluster alongside FAPs/sCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (Lamc2/Itgb4)
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Lamc2"
$ws.Cells.Item(2, 3).Value = "Itgb4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2.0
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 0.2121133333333333
$ws.Cells.Item(2, 8).Value = 0.6363399999999999
$ws.Cells.Item(2, 9).Value = 0.04296779043029777
$ws.Cells.Item(2, 10).Value = 0.04296779043029776
$ws.Cells.Item(2, 11).Value = 2.0
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 4.187598666666667
$ws.Cells.Item(2, 14).Value = 12.562796
$ws.Cells.Item(2, 15).Value = 0.4284941686600627
$ws.Cells.Item(2, 16).Value = 0.4284941686600626
$ws.Cells.Item(2, 17).Value = 0.8882455118488888
$ws.Cells.Item(2, 18).Value = 7.994209606639999
$ws.Cells.Item(2, 19).Value = 0.01841144763959024
$ws.Cells.Item(2, 20).Value = 0.01841144763959023

# Row 3: ECs -> FAPs (Lamc2/Itgb4)
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Lamc2"
$ws.Cells.Item(3, 3).Value = "Itgb4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2.0
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 0.2121133333333333
$ws.Cells.Item(3, 8).Value = 0.6363399999999999
$ws.Cells.Item(3, 9).Value = 0.04296779043029777
$ws.Cells.Item(3, 10).Value = 0.04296779043029776
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 1.204585666666667
$ws.Cells.Item(3, 14).Value = 3.613757
$ws.Cells.Item(3, 15).Value = 0.1232586918910792
$ws.Cells.Item(3, 16).Value = 0.1232586918910792
$ws.Cells.Item(3, 17).Value = 0.2555086810422222
$ws.Cells.Item(3, 18).Value = 2.29957812938
$ws.Cells.Item(3, 19).Value = 0.005296153641888534
$ws.Cells.Item(3, 20).Value = 0.005296153641888532

# Row 4: ECs -> sCs (Lamc2/Itgb4)
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Lamc2"
$ws.Cells.Item(4, 3).Value = "Itgb4"
$ws.Cells.Item(4, 4).Value = "sCs"
$ws.Cells.Item(4, 5).Value = 2.0
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 0.2121133333333333
$ws.Cells.Item(4, 8).Value = 0.6363399999999999
$ws.Cells.Item(4, 9).Value = 0.04296779043029777
$ws.Cells.Item(4, 10).Value = 0.04296779043029776
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 12).Value = 1.0
$ws.Cells.Item(4, 13).Value = 4.380641
$ws.Cells.Item(4, 14).Value = 13.141923
$ws.Cells.Item(4, 15).Value = 0.4482471394488581
$ws.Cells.Item(4, 16).Value = 0.4482471394488581
$ws.Cells.Item(4, 17).Value = 0.9291923646466664
$ws.Cells.Item(4, 18).Value = 8.362731281819999
$ws.Cells.Item(4, 19).Value = 0.01926018914881899
$ws.Cells.Item(4, 20).Value = 0.01926018914881899

# Row 5: FAPs -> ECs (Lamc2/Itgb4)
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Lamc2"
$ws.Cells.Item(5, 3).Value = "Itgb4"
$ws.Cells.Item(5, 4).Value = "ECs"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 3.865410333333333
$ws.Cells.Item(5, 8).Value = 11.596231
$ws.Cells.Item(5, 9).Value = 0.7830160344930734
$ws.Cells.Item(5, 10).Value = 0.7830160344930733
$ws.Cells.Item(5, 11).Value = 2.0
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 4.187598666666667
$ws.Cells.Item(5, 14).Value = 12.562796
$ws.Cells.Item(5, 15).Value = 0.4284941686600627
$ws.Cells.Item(5, 16).Value = 0.4284941686600626
$ws.Cells.Item(5, 17).Value = 16.18678715798622
$ws.Cells.Item(5, 18).Value = 145.681084421876
$ws.Cells.Item(5, 19).Value = 0.3355178047476085
$ws.Cells.Item(5, 20).Value = 0.3355178047476084

# Row 6: FAPs -> FAPs (Lamc2/Itgb4)
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Lamc2"
$ws.Cells.Item(6, 3).Value = "Itgb4"
$ws.Cells.Item(6, 4).Value = "FAPs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 3.865410333333333
$ws.Cells.Item(6, 8).Value = 11.596231
$ws.Cells.Item(6, 9).Value = 0.7830160344930734
$ws.Cells.Item(6, 10).Value = 0.7830160344930733
$ws.Cells.Item(6, 11).Value = 3.0
$ws.Cells.Item(6, 12).Value = 1.0
$ws.Cells.Item(6, 13).Value = 1.204585666666667
$ws.Cells.Item(6, 14).Value = 3.613757
$ws.Cells.Item(6, 15).Value = 0.1232586918910792
$ws.Cells.Item(6, 16).Value = 0.1232586918910792
$ws.Cells.Item(6, 17).Value = 4.656217883318556
$ws.Cells.Item(6, 18).Value = 41.905960949867
$ws.Cells.Item(6, 19).Value = 0.09651353214135638
$ws.Cells.Item(6, 20).Value = 0.09651353214135636

# Row 7: FAPs -> sCs (Lamc2/Itgb4)
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Lamc2"
$ws.Cells.Item(7, 3).Value = "Itgb4"
$ws.Cells.Item(7, 4).Value = "sCs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 3.865410333333333
$ws.Cells.Item(7, 8).Value = 11.596231
$ws.Cells.Item(7, 9).Value = 0.7830160344930734
$ws.Cells.Item(7, 10).Value = 0.7830160344930733
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 4.380641
$ws.Cells.Item(7, 14).Value = 13.141923
$ws.Cells.Item(7, 15).Value = 0.4482471394488581
$ws.Cells.Item(7, 16).Value = 0.4482471394488581
$ws.Cells.Item(7, 17).Value = 16.93297498802367
$ws.Cells.Item(7, 18).Value = 152.396774892213
$ws.Cells.Item(7, 19).Value = 0.3509846976041086
$ws.Cells.Item(7, 20).Value = 0.3509846976041085

# Row 8: sCs -> ECs (Lamc2/Itgb4)
$ws.Cells.Item(8, 1).Value = "sCs"
$ws.Cells.Item(8, 2).Value = "Lamc2"
$ws.Cells.Item(8, 3).Value = "Itgb4"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 0.8590423333333334
$ws.Cells.Item(8, 8).Value = 2.577127
$ws.Cells.Item(8, 9).Value = 0.1740161750766289
$ws.Cells.Item(8, 10).Value = 0.1740161750766288
$ws.Cells.Item(8, 11).Value = 2.0
$ws.Cells.Item(8, 12).Value = 0.6666666666666666
$ws.Cells.Item(8, 13).Value = 4.187598666666667
$ws.Cells.Item(8, 14).Value = 12.562796
$ws.Cells.Item(8, 15).Value = 0.4284941686600627
$ws.Cells.Item(8, 16).Value = 0.4284941686600626
$ws.Cells.Item(8, 17).Value = 3.597324529676889
$ws.Cells.Item(8, 18).Value = 32.375920767092
$ws.Cells.Item(8, 19).Value = 0.07456491627286402
$ws.Cells.Item(8, 20).Value = 0.07456491627286398

# Row 9: sCs -> FAPs (Lamc2/Itgb4)
$ws.Cells.Item(9, 1).Value = "sCs"
$ws.Cells.Item(9, 2).Value = "Lamc2"
$ws.Cells.Item(9, 3).Value = "Itgb4"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 0.8590423333333334
$ws.Cells.Item(9, 8).Value = 2.577127
$ws.Cells.Item(9, 9).Value = 0.1740161750766289
$ws.Cells.Item(9, 10).Value = 0.1740161750766288
$ws.Cells.Item(9, 11).Value = 3.0
$ws.Cells.Item(9, 12).Value = 1.0
$ws.Cells.Item(9, 13).Value = 1.204585666666667
$ws.Cells.Item(9, 14).Value = 3.613757
$ws.Cells.Item(9, 15).Value = 0.1232586918910792
$ws.Cells.Item(9, 16).Value = 0.1232586918910792
$ws.Cells.Item(9, 17).Value = 1.034790081793222
$ws.Cells.Item(9, 18).Value = 9.313110736139
$ws.Cells.Item(9, 19).Value = 0.02144900610783429
$ws.Cells.Item(9, 20).Value = 0.02144900610783429

# Row 10: sCs -> sCs (Lamc2/Itgb4)
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Lamc2"
$ws.Cells.Item(10, 3).Value = "Itgb4"
$ws.Cells.Item(10, 4).Value = "sCs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 0.8590423333333334
$ws.Cells.Item(10, 8).Value = 2.577127
$ws.Cells.Item(10, 9).Value = 0.1740161750766289
$ws.Cells.Item(10, 10).Value = 0.1740161750766288
$ws.Cells.Item(10, 11).Value = 3.0
$ws.Cells.Item(10, 12).Value = 1.0
$ws.Cells.Item(10, 13).Value = 4.380641
$ws.Cells.Item(10, 14).Value = 13.141923
$ws.Cells.Item(10, 15).Value = 0.4482471394488581
$ws.Cells.Item(10, 16).Value = 0.4482471394488581
$ws.Cells.Item(10, 17).Value = 3.763156066135667
$ws.Cells.Item(10, 18).Value = 33.868404595221
$ws.Cells.Item(10, 19).Value = 0.07800225269593056
$ws.Cells.Item(10, 20).Value = 0.07800225269593054

